# template_barang.xlsx: drop the extra leading zero in the generated
# barang_kode (item code) values, e.g. "BRG0011" -> "BRG011", and leave
# the selection on the last edited cell (B6) as Excel would after typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BRG011"
$ws.Range("B3").Value = "BRG012"
$ws.Range("B4").Value = "BRG013"
$ws.Range("B5").Value = "BRG014"
$ws.Range("B6").Value = "BRG015"

$ws.Range("B6").Select()
